$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 ("player16"): the ACTUAL OUTPUT changes from the predicted
# "player16 remains on the list of online players" to what really happened:
# "player16 is removed from the list of online players". The old text was
# predicted (non-bold) and this one also ends up non-bold (border kept).
$ws.Range("C18").Font.Bold = $true
$ws.Range("C18").Value = "player16 is removed from the list of online players"
$ws.Range("C18").Font.Bold = $false

# --- Row 19 was a blank spacer row between the "player16" block and the
# "player17" block; remove it so player17..player23 rows shift up by one.
$ws.Rows(19).Delete()

# --- Put the final selection where the author left off.
$ws.Range("A26").Select()
